$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(45736,45736,45736,45736,45736,45736,45736,45736,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45737,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45738,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45739,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45740,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45741,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45742,45743,45743,45743,45743,45743,45743,45743,45743,45743,45743,45743,45743,45743,45743,45743,45743,45743)
$colB = @(17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17)
$colC = @(0.011,0.629,0.054,0.011,0.011,0.011,0.011,0.011,0.011,0.011,0.011,0.011,0.011,0.011,0.021,0.583,1.581,2.647,3.366,3.95,4.181,4.156,3.819,3.032,1.72,0.662,0.054,0.011,0.011,0.011,0.011,0.011,0.011,0.011,0.011,0.012,0.012,0.012,0.013,0.443,1.451,2.163,2.484,2.341,2.148,1.694,1.368,1.115,0.77,0.205,0.014,0.013,0.013,0.013,0.013,0.013,0.013,0.013,0.013,0.012,0.012,0.012,0.013,0.195,0.766,1.282,1.513,1.622,1.529,1.41,1.114,0.842,0.596,0.132,0.014,0.013,0.013,0.013,0.013,0.013,0.013,0.013,0.013,0.013,0.013,0.013,0.013,0.128,0.351,0.757,0.871,1.219,1.482,1.427,1.349,1.208,0.783,0.252,0.023,0.013,0.013,0.013,0.013,0.013,0.012,0.012,0.012,0.012,0.012,0.012,0.022,0.45,1.376,2.179,2.918,3.113,3.154,3.055,2.49,1.52,0.783,0.2,0.017,0.013,0.013,0.013,0.013,0.013,0.013,0.013,0.013,0.013,0.012,0.012,0.013,0.247,0.922,1.752,2.607,3.153,3.349,3.327,3.014,2.272,1.427,0.439,0.037,0.012,0.012,0.012,0.012,0.012,0.012,0.012,0.013,0.013,0.013,0.013,0.014,0.116,0.312,0.639,0.8080000000000001,0.996,1.02,0.971,0.801,0.734,0.369)
$colD = @("20.03.202517","20.03.202518","20.03.202519","20.03.202520","20.03.202521","20.03.202522","20.03.202523","20.03.202524","21.03.20251","21.03.20252","21.03.20253","21.03.20254","21.03.20255","21.03.20256","21.03.20257","21.03.20258","21.03.20259","21.03.202510","21.03.202511","21.03.202512","21.03.202513","21.03.202514","21.03.202515","21.03.202516","21.03.202517","21.03.202518","21.03.202519","21.03.202520","21.03.202521","21.03.202522","21.03.202523","21.03.202524","22.03.20251","22.03.20252","22.03.20253","22.03.20254","22.03.20255","22.03.20256","22.03.20257","22.03.20258","22.03.20259","22.03.202510","22.03.202511","22.03.202512","22.03.202513","22.03.202514","22.03.202515","22.03.202516","22.03.202517","22.03.202518","22.03.202519","22.03.202520","22.03.202521","22.03.202522","22.03.202523","22.03.202524","23.03.20251","23.03.20252","23.03.20253","23.03.20254","23.03.20255","23.03.20256","23.03.20257","23.03.20258","23.03.20259","23.03.202510","23.03.202511","23.03.202512","23.03.202513","23.03.202514","23.03.202515","23.03.202516","23.03.202517","23.03.202518","23.03.202519","23.03.202520","23.03.202521","23.03.202522","23.03.202523","23.03.202524","24.03.20251","24.03.20252","24.03.20253","24.03.20254","24.03.20255","24.03.20256","24.03.20257","24.03.20258","24.03.20259","24.03.202510","24.03.202511","24.03.202512","24.03.202513","24.03.202514","24.03.202515","24.03.202516","24.03.202517","24.03.202518","24.03.202519","24.03.202520","24.03.202521","24.03.202522","24.03.202523","24.03.202524","25.03.20251","25.03.20252","25.03.20253","25.03.20254","25.03.20255","25.03.20256","25.03.20257","25.03.20258","25.03.20259","25.03.202510","25.03.202511","25.03.202512","25.03.202513","25.03.202514","25.03.202515","25.03.202516","25.03.202517","25.03.202518","25.03.202519","25.03.202520","25.03.202521","25.03.202522","25.03.202523","25.03.202524","26.03.20251","26.03.20252","26.03.20253","26.03.20254","26.03.20255","26.03.20256","26.03.20257","26.03.20258","26.03.20259","26.03.202510","26.03.202511","26.03.202512","26.03.202513","26.03.202514","26.03.202515","26.03.202516","26.03.202517","26.03.202518","26.03.202519","26.03.202520","26.03.202521","26.03.202522","26.03.202523","26.03.202524","27.03.20251","27.03.20252","27.03.20253","27.03.20254","27.03.20255","27.03.20256","27.03.20257","27.03.20258","27.03.20259","27.03.202510","27.03.202511","27.03.202512","27.03.202513","27.03.202514","27.03.202515","27.03.202516","27.03.202517")

for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
    $ws.Cells.Item($r, 4).Value = $colD[$i]
}
